# Point d'etape projet WIOD TIVA - homogeneisation des echelles et ajouts de texte explicatif
#
# This script is applied to the *original* 10-slide deck. Shape/slide indices
# below refer to the ORIGINAL layout; we perform in-place text tidy-ups first
# (run merges that don't change the visible text, just collapse split runs),
# then insert the two brand-new slides last so earlier indices stay stable.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 2 ("Donnees" table): collapse a few split runs into single runs.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tbl = $s2.Shapes.Item(2).Table

# Row "Secteurs" / col WIOD -> "Agriculture, biens, services" (was 3 runs)
$tbl.Cell(4,1).Shape.TextFrame.TextRange.Paragraphs(2,1).Text = "Agriculture, biens, services"

# Row "Import proportionality assumption" col 1 (was 2 runs)
$tbl.Cell(5,1).Shape.TextFrame.TextRange.Paragraphs(1,1).Text = "Import proportionality assumption"

# Row col 2, second paragraph: "uses bilateral trade statistics to derive import shares for three end-use categories " (was 2 runs; 3rd run with different formatting stays separate)
$tbl.Cell(5,2).Shape.TextFrame.TextRange.Paragraphs(2,1).Text = "uses bilateral trade statistics to derive import shares for three end-use categories (intermediate use, ﬁnal consumption and investment)"

# Row col 3, second paragraph: "assumes that the share of imports..." (was 2 runs)
$tbl.Cell(5,3).Shape.TextFrame.TextRange.Paragraphs(2,1).Text = "assumes that the share of imports in any product consumed directly as intermediate consumption or final demand is the same for all users"

# ---------------------------------------------------------------------------
# 2) Slide 5 title: merge 3 runs into one (no visible text change)
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Élasticité prix conso pays hors ZE à une appréciation EUR"

# ---------------------------------------------------------------------------
# 3) Slide 9 title: merge 2 runs, appending "USD" into the same sentence
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Élasticité prix conso ZE à une appréciation USD"

# ---------------------------------------------------------------------------
# 4) Slide 10 title: merge 2 runs, appending "UK" into the same sentence
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Élasticité prix conso ZE à une appréciation UK"

# ---------------------------------------------------------------------------
# 5) Insert a new slide at position 2: "Travail effectué"
# ---------------------------------------------------------------------------
$titleContentLayout = $p.SlideMaster.CustomLayouts.Item(2)
$sTravail = $p.Slides.AddSlide(2, $titleContentLayout)
$sTravail.Shapes.Item(1).TextFrame.TextRange.Text = "Travail effectué"

$bodyTravail = $sTravail.Shapes.Item(2)
$bodyTravail.Left = 457200
$bodyTravail.Top = 1268760
$bodyTravail.Width = 8229600
$bodyTravail.Height = 4857403

$trTravail = $bodyTravail.TextFrame.TextRange
$trTravail.Text = "Reproduction du travail fourni dans le WP OFCE à une nouvelle base de données (WIOD) et prise en compte des mises à jour de la base TIVA`rImpact de chocs de change sur:`rPrix de production`rPrix d’exportations`rPrix de consommation`r"
$trTravail.Paragraphs(3,1).IndentLevel = 2
$trTravail.Paragraphs(4,1).IndentLevel = 2
$trTravail.Paragraphs(5,1).IndentLevel = 2

# ---------------------------------------------------------------------------
# 6) Insert a new slide at the end: "Suite du projet"
# ---------------------------------------------------------------------------
$sSuite = $p.Slides.AddSlide($p.Slides.Count + 1, $titleContentLayout)
$sSuite.Shapes.Item(1).TextFrame.TextRange.Text = "Suite du projet"

$bodySuite = $sSuite.Shapes.Item(2)
$trSuite = $bodySuite.TextFrame.TextRange
$trSuite.Text = "Interprétation des résultats`rApprofondir l’analyse sur la part des inputs importés `rImpact d’un choc sur le prix du pétrole sur les déflateurs"
